# Update factsheets with text edits from COMM
#
# Summary of changes:
#  - "Overall" sheet: A2 (filer count) becomes text "4,025" instead of numeric 4025
#  - "County" sheet: all numeric filer-count cells in column B become text; the
#    Wyandot County row (row 89, all zeros) gets reformatted to percent/currency
#    strings; a new "Total" row (row 90) is appended with the overall totals.
#  - "Congressional District", "Size", "Subsector" sheets: all numeric filer-count
#    cells in column B become text (including their "Total" rows).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overall"
# ---------------------------------------------------------------------------
$wsOverall = $wb.Worksheets.Item("Overall")
$rng = $wsOverall.Range("A2")
$rng.NumberFormat = "@"
$rng.Value = "4,025"

# ---------------------------------------------------------------------------
# Sheet "County"
# ---------------------------------------------------------------------------
$wsCounty = $wb.Worksheets.Item("County")

$countyCounts = @{
    2 = "2"; 3 = "42"; 4 = "18"; 5 = "30"; 6 = "38"; 7 = "15"; 8 = "11"; 9 = "7";
    10 = "61"; 11 = "5"; 12 = "12"; 13 = "42"; 14 = "24"; 15 = "10"; 16 = "30";
    17 = "12"; 18 = "12"; 19 = "550"; 20 = "18"; 21 = "13"; 22 = "50"; 23 = "25";
    24 = "33"; 25 = "9"; 26 = "726"; 27 = "12"; 28 = "14"; 29 = "25"; 30 = "39";
    31 = "17"; 32 = "383"; 33 = "30"; 34 = "12"; 35 = "2"; 36 = "13"; 37 = "13";
    38 = "7"; 39 = "10"; 40 = "11"; 41 = "15"; 42 = "29"; 43 = "27"; 44 = "37";
    45 = "6"; 46 = "45"; 47 = "14"; 48 = "58"; 49 = "188"; 50 = "9"; 51 = "71";
    52 = "68"; 53 = "29"; 54 = "2"; 55 = "13"; 56 = "30"; 57 = "5"; 58 = "171";
    59 = "5"; 60 = "4"; 61 = "25"; 62 = "5"; 63 = "31"; 64 = "3"; 65 = "7";
    66 = "10"; 67 = "7"; 68 = "30"; 69 = "7"; 70 = "10"; 71 = "46"; 72 = "24";
    73 = "16"; 74 = "26"; 75 = "25"; 76 = "14"; 77 = "119"; 78 = "192"; 79 = "30";
    80 = "20"; 81 = "14"; 82 = "18"; 83 = "3"; 84 = "34"; 85 = "20"; 86 = "36";
    87 = "10"; 88 = "34"
}

foreach ($row in $countyCounts.Keys) {
    $cellRef = "B" + $row
    $rng = $wsCounty.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $countyCounts[$row]
}

# Row 89 (Wyandot County) - reformat the all-zero stats as percent/currency text
$rng = $wsCounty.Range("B89")
$rng.NumberFormat = "@"
$rng.Value = "0.00%"

$rng = $wsCounty.Range("C89")
$rng.NumberFormat = "@"
$rng.Value = "`$0"

$rng = $wsCounty.Range("D89")
$rng.NumberFormat = "@"
$rng.Value = "0.00%"

$rng = $wsCounty.Range("E89")
$rng.NumberFormat = "@"
$rng.Value = "0.00%"

$rng = $wsCounty.Range("F89")
$rng.NumberFormat = "@"
$rng.Value = "0.00%"

# New row 90 - overall Total row
$rng = $wsCounty.Range("A90")
$rng.NumberFormat = "@"
$rng.Value = "Total"

$rng = $wsCounty.Range("B90")
$rng.NumberFormat = "@"
$rng.Value = "4,025"

$rng = $wsCounty.Range("C90")
$rng.NumberFormat = "@"
$rng.Value = "`$7,999,510,499"

$rng = $wsCounty.Range("D90")
$rng.NumberFormat = "@"
$rng.Value = "7.17%"

$rng = $wsCounty.Range("E90")
$rng.NumberFormat = "@"
$rng.Value = "-17.50%"

$rng = $wsCounty.Range("F90")
$rng.NumberFormat = "@"
$rng.Value = "71.45%"

# ---------------------------------------------------------------------------
# Sheet "Congressional District"
# ---------------------------------------------------------------------------
$wsCd = $wb.Worksheets.Item("Congressional District")

$cdCounts = @{
    2 = "372"; 3 = "250"; 4 = "463"; 5 = "220"; 6 = "288"; 7 = "152"; 8 = "222";
    9 = "170"; 10 = "565"; 11 = "295"; 12 = "218"; 13 = "220"; 14 = "153";
    15 = "133"; 16 = "304"
}
foreach ($row in $cdCounts.Keys) {
    $cellRef = "B" + $row
    $rng = $wsCd.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $cdCounts[$row]
}

$rng = $wsCd.Range("B17")
$rng.NumberFormat = "@"
$rng.Value = "4,025"

# ---------------------------------------------------------------------------
# Sheet "Size"
# ---------------------------------------------------------------------------
$wsSize = $wb.Worksheets.Item("Size")

$sizeCounts = @{
    2 = "1,280"; 3 = "1,057"; 4 = "644"; 5 = "313"; 6 = "541"; 7 = "190"
}
foreach ($row in $sizeCounts.Keys) {
    $cellRef = "B" + $row
    $rng = $wsSize.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $sizeCounts[$row]
}

$rng = $wsSize.Range("B8")
$rng.NumberFormat = "@"
$rng.Value = "4,025"

# ---------------------------------------------------------------------------
# Sheet "Subsector"
# ---------------------------------------------------------------------------
$wsSub = $wb.Worksheets.Item("Subsector")

$subCounts = @{
    2 = "341"; 3 = "458"; 4 = "126"; 5 = "365"; 6 = "41"; 7 = "1,332"; 8 = "16";
    9 = "1"; 10 = "341"; 11 = "114"; 12 = "837"; 13 = "53"
}
foreach ($row in $subCounts.Keys) {
    $cellRef = "B" + $row
    $rng = $wsSub.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $subCounts[$row]
}

$rng = $wsSub.Range("B14")
$rng.NumberFormat = "@"
$rng.Value = "4,025"
